$wb = $excel.ActiveWorkbook

# ALC row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3061
$ws.Range("I19").Value = 3380.1667
$ws.Range("K19").Value = 3380.1667
$ws.Range("M19").Value = -3205.1667

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 71758.13
$ws.Range("I62").Value = 95261.82000000001
$ws.Range("K62").Value = 95261.82000000001
$ws.Range("M62").Value = -94637.82000000001

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 71758.13
$ws.Range("I65").Value = 95261.82000000001
$ws.Range("K65").Value = 476309.1
$ws.Range("M65").Value = -473189.1

# ALC row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4101.6665
$ws.Range("I74").Value = 3065.875
$ws.Range("J74").Value = 5285.4287
$ws.Range("K74").Value = 3065.875
$ws.Range("L74").Value = 5285.4287
$ws.Range("M74").Value = -2129.875
$ws.Range("N74").Value = -7157.4287

# ALC row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 4101.6665
$ws.Range("I77").Value = 3065.875
$ws.Range("J77").Value = 5285.4287
$ws.Range("K77").Value = 15329.375
$ws.Range("L77").Value = 26427.1435
$ws.Range("M77").Value = -10649.375
$ws.Range("N77").Value = -35787.14350000001

# ALC row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2508.2222
$ws.Range("I86").Value = 2512.3333
$ws.Range("J86").Value = 2500
$ws.Range("K86").Value = 2512.3333
$ws.Range("L86").Value = 2500
$ws.Range("M86").Value = -1389.3333
$ws.Range("N86").Value = -4746

# ALC row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 2508.2222
$ws.Range("I89").Value = 2512.3333
$ws.Range("J89").Value = 2500
$ws.Range("K89").Value = 12561.6665
$ws.Range("L89").Value = 12500
$ws.Range("M89").Value = -6945.666499999999
$ws.Range("N89").Value = -23732

# ALC row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 795.7273
$ws.Range("I92").Value = 775.3
$ws.Range("J92").Value = 1000
$ws.Range("K92").Value = 775.3
$ws.Range("L92").Value = 1000
$ws.Range("M92").Value = 472.7
$ws.Range("N92").Value = -3496

# ALC row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 22798.8
$ws.Range("I106").Value = 2499.5
$ws.Range("K106").Value = 2499.5
$ws.Range("M106").Value = -1868.5

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2379.375
$ws.Range("J2").Value = 2437.25
$ws.Range("L2").Value = 2437.25
$ws.Range("N2").Value = -2663.25

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 38541.605
$ws.Range("I32").Value = 21186.76
$ws.Range("K32").Value = 21186.76
$ws.Range("M32").Value = -20899.76

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 419444.88
$ws.Range("I45").Value = 558076.5
$ws.Range("J45").Value = 3550
$ws.Range("K45").Value = 558076.5
$ws.Range("L45").Value = 3550
$ws.Range("M45").Value = -557699.5
$ws.Range("N45").Value = -4304

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1750
$ws.Range("J61").Value = 2000
$ws.Range("L61").Value = 2000
$ws.Range("N61").Value = -2424

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1925.6
$ws.Range("I63").Value = 1937.6923
$ws.Range("K63").Value = 1937.6923
$ws.Range("M63").Value = -1251.6923

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 1925.6
$ws.Range("I66").Value = 1937.6923
$ws.Range("K66").Value = 9688.461499999999
$ws.Range("M66").Value = -6256.461499999999

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2379.375
$ws.Range("J116").Value = 2437.25
$ws.Range("L116").Value = 2437.25
$ws.Range("N116").Value = -7025.25

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1887.8846
$ws.Range("I122").Value = 1836.3334
$ws.Range("J122").Value = 2506.5
$ws.Range("K122").Value = 5509.0002
$ws.Range("L122").Value = 7519.5
$ws.Range("M122").Value = -3059.0002
$ws.Range("N122").Value = -12419.5

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 25767.68
$ws.Range("I132").Value = 28767.455
$ws.Range("K132").Value = 86302.36500000001
$ws.Range("M132").Value = -83772.36500000001

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1750
$ws.Range("J136").Value = 2000
$ws.Range("L136").Value = 6000
$ws.Range("N136").Value = -11100

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2379.375
$ws.Range("J3").Value = 2437.25
$ws.Range("L3").Value = 2437.25
$ws.Range("N3").Value = -2665.25

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2262.7646
$ws.Range("I86").Value = 1881.6154
$ws.Range("J86").Value = 3501.5
$ws.Range("K86").Value = 1881.6154
$ws.Range("L86").Value = 3501.5
$ws.Range("M86").Value = -758.6153999999999
$ws.Range("N86").Value = -5747.5

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2262.7646
$ws.Range("I89").Value = 1881.6154
$ws.Range("J89").Value = 3501.5
$ws.Range("K89").Value = 9408.076999999999
$ws.Range("L89").Value = 17507.5
$ws.Range("M89").Value = -3792.076999999999
$ws.Range("N89").Value = -28739.5

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 468.08334
$ws.Range("I94").Value = 468.08334
$ws.Range("K94").Value = 468.08334
$ws.Range("M94").Value = -17.08334000000002

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1482.742
$ws.Range("I31").Value = 1498.9615
$ws.Range("J31").Value = 1398.4
$ws.Range("K31").Value = 1498.9615
$ws.Range("L31").Value = 1398.4
$ws.Range("M31").Value = -1203.9615
$ws.Range("N31").Value = -1988.4

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1482.742
$ws.Range("I34").Value = 1498.9615
$ws.Range("J34").Value = 1398.4
$ws.Range("K34").Value = 1498.9615
$ws.Range("L34").Value = 1398.4
$ws.Range("M34").Value = -1296.9615
$ws.Range("N34").Value = -1802.4

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2213.4348
$ws.Range("I122").Value = 2603.4614
$ws.Range("J122").Value = 1706.4
$ws.Range("K122").Value = 7810.3842
$ws.Range("L122").Value = 5119.200000000001
$ws.Range("M122").Value = -5360.3842
$ws.Range("N122").Value = -10019.2

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2558.524
$ws.Range("I132").Value = 2353.4211
$ws.Range("K132").Value = 7060.263300000001
$ws.Range("M132").Value = -4530.263300000001

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2361.85
$ws.Range("I134").Value = 2334
$ws.Range("J134").Value = 2473.25
$ws.Range("K134").Value = 7002
$ws.Range("L134").Value = 7419.75
$ws.Range("M134").Value = -4467
$ws.Range("N134").Value = -12489.75

# CUL row 2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 23.75
$ws.Range("J2").Value = 23.75
$ws.Range("L2").Value = 142.5
$ws.Range("N2").Value = -368.5

# CUL row 74
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 4670.6665
$ws.Range("I74").Value = 4998.5
$ws.Range("J74").Value = 4015
$ws.Range("K74").Value = 14995.5
$ws.Range("L74").Value = 12045
$ws.Range("M74").Value = -13934.5
$ws.Range("N74").Value = -14167

# CUL row 77
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H77").Value = 4670.6665
$ws.Range("I77").Value = 4998.5
$ws.Range("J77").Value = 4015
$ws.Range("K77").Value = 44986.5
$ws.Range("L77").Value = 36135
$ws.Range("M77").Value = -39682.5
$ws.Range("N77").Value = -46743

# CUL row 87
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 5498.5
$ws.Range("I87").Value = 4998.3335
$ws.Range("K87").Value = 14995.0005
$ws.Range("M87").Value = -13747.0005

# CUL row 90
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 5498.5
$ws.Range("I90").Value = 4998.3335
$ws.Range("K90").Value = 44985.0015
$ws.Range("M90").Value = -38745.0015

# CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 780.9375
$ws.Range("J107").Value = 883.3913
$ws.Range("L107").Value = 2650.1739
$ws.Range("N107").Value = -6490.1739

# CUL row 115
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H115").Value = 2228.5715
$ws.Range("I115").Value = 1200
$ws.Range("K115").Value = 3600
$ws.Range("M115").Value = -2425

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2992.25
$ws.Range("I126").Value = 2836.6924
$ws.Range("J126").Value = 3666.3333
$ws.Range("K126").Value = 8510.0772
$ws.Range("L126").Value = 10998.9999
$ws.Range("M126").Value = -6040.0772
$ws.Range("N126").Value = -15938.9999

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 700
$ws.Range("J61").Value = 700
$ws.Range("L61").Value = 700
$ws.Range("N61").Value = -1104

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 700
$ws.Range("J113").Value = 700
$ws.Range("L113").Value = 700
$ws.Range("N113").Value = -5040

# WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7205.8096
$ws.Range("I62").Value = 7675.933
$ws.Range("K62").Value = 7675.933
$ws.Range("M62").Value = -7051.933

# WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 7205.8096
$ws.Range("I65").Value = 7675.933
$ws.Range("K65").Value = 38379.665
$ws.Range("M65").Value = -35259.665

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3339.125
$ws.Range("I81").Value = 3339.125
$ws.Range("K81").Value = 6678.25
$ws.Range("M81").Value = -5617.25

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 3339.125
$ws.Range("I84").Value = 3339.125
$ws.Range("K84").Value = 33391.25
$ws.Range("M84").Value = -28087.25
